$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing row 2 timestamp ("Data Criação") ---
$ws.Range("N2").Value = "2025-12-09 09:46:42"

# --- Prime row 3 with row 2's formatting (borders/alignment/row height) ---
$ws.Range("A2:N2").Copy()
$ws.Range("A3:N3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(3).RowHeight = 20

# --- Fill in row 3 values ---
$ws.Range("A3").Value = "001_Recursos_proprios"

# Text-like values that would otherwise be mis-typed as numbers/dates by
# COM's automatic type inference: force them to stay text via NumberFormat.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "000"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "0000000000000000"

$ws.Range("D3").Value = "carlos"

# E3 / F3 are intentionally left blank (Programa / Campus not informed)
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""

$ws.Range("G3").Value = "mestrado"
$ws.Range("H3").Value = "R$ 2.100.00"

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "2025-04-05"

$ws.Range("J3").Value = "2025-12"

# K3 ("Parcela") is a real number
$ws.Range("K3").Value = 9

$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = "444"

# M3 ("Observações") left blank
$ws.Range("M3").Value = ""

$ws.Range("N3").Value = "2025-12-09 23:58:21"

# --- Re-apply row 2's formats on top of the new values so every cell in
#     row 3 ends up on the exact same style (border/alignment/number
#     format) as its row-2 counterpart, instead of a freshly minted one. ---
$ws.Range("A2:N2").Copy()
$ws.Range("A3:N3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(3).RowHeight = 20
